# Apply the "cryptos list" data refresh described in the commit:
# "Updated cryptos list on Tue Nov  5 14:38:45 UTC 2024 with GitHub Actions"
#
# Updates the Price (D) / Volume(1h) (E) columns for most rows, and swaps the
# ranking of "WrappedeETH" (row 26) and "Aptos" (row 27), which also moves
# along their Coin (B) / Link (C) / Price (D) / Volume (E) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, even when it looks like a number
# (e.g. "567.47"), so Excel does not silently convert the cell to a numeric
# type. The cell keeps its original (default) style afterwards.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '69.679.40'
$ws.Range('E2').Value = '  +1.77%  '
$ws.Range('D3').Value = '2.467.06'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  -0.20%  '
Set-TextValue ($ws.Range('D5')) '567.47'
$ws.Range('E5').Value = '  +1.67%  '
Set-TextValue ($ws.Range('D6')) '166.38'
$ws.Range('E6').Value = '  +2.14%  '
Set-TextValue ($ws.Range('D7')) '1.00'
Set-TextValue ($ws.Range('D8')) '0.512'
$ws.Range('E8').Value = '  +0.41%  '
Set-TextValue ($ws.Range('D9')) '0.177'
$ws.Range('E9').Value = '  +12.89%  '
Set-TextValue ($ws.Range('D10')) '0.162'
$ws.Range('E10').Value = '  -1.55%  '
Set-TextValue ($ws.Range('D11')) '0.335'
$ws.Range('E11').Value = '  +1.94%  '
Set-TextValue ($ws.Range('D12')) '4.68'
$ws.Range('E12').Value = '  -2.95%  '
Set-TextValue ($ws.Range('D13')) '0.0000183'
$ws.Range('E13').Value = '  +8.03%  '
$ws.Range('D14').Value = '69.501.53'
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').Value = '2.912.68'
$ws.Range('E15').Value = '  +0.24%  '
Set-TextValue ($ws.Range('D16')) '23.81'
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('D17').Value = '2.451.44'
$ws.Range('E17').Value = '  -0.48%  '
Set-TextValue ($ws.Range('D18')) '10.81'
$ws.Range('E18').Value = '  +2.98%  '
Set-TextValue ($ws.Range('D19')) '342.70'
$ws.Range('E19').Value = '  +1.68%  '
Set-TextValue ($ws.Range('D20')) '7.16'
$ws.Range('E20').Value = '  +4.17%  '
Set-TextValue ($ws.Range('D21')) '3.90'
$ws.Range('E21').Value = '  +2.90%  '
Set-TextValue ($ws.Range('D22')) '2.01'
$ws.Range('E22').Value = '  +7.04%  '
$ws.Range('E23').Value = '  +0.00%  '
Set-TextValue ($ws.Range('D24')) '66.19'
$ws.Range('E24').Value = '  -0.70%  '
Set-TextValue ($ws.Range('D25')) '3.88'
$ws.Range('E25').Value = '  +6.41%  '
$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue ($ws.Range('D26')) '8.56'
$ws.Range('E26').Value = '  +5.41%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.573.68'
$ws.Range('E27').Value = '  -0.34%  '
Set-TextValue ($ws.Range('D28')) '0.986'
$ws.Range('E28').Value = '  -1.42%  '
$ws.Range('D29').Value = '0.0₃0850'
$ws.Range('E29').Value = '  +3.99%  '
Set-TextValue ($ws.Range('D30')) '7.32'
$ws.Range('E30').Value = '  +1.89%  '
Set-TextValue ($ws.Range('D31')) '1.25'
$ws.Range('E31').Value = '  +9.75%  '
Set-TextValue ($ws.Range('D32')) '448.86'
$ws.Range('E32').Value = '  +5.74%  '
Set-TextValue ($ws.Range('D33')) '0.999'
$ws.Range('E33').Value = '  -0.07%  '
Set-TextValue ($ws.Range('D34')) '1.63'
$ws.Range('E34').Value = '  +0.82%  '
Set-TextValue ($ws.Range('D35')) '160.20'
$ws.Range('E35').Value = '  -0.55%  '
Set-TextValue ($ws.Range('D36')) '19.04'
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('E37').Value = '  -0.03%  '
Set-TextValue ($ws.Range('D38')) '0.110'
$ws.Range('E38').Value = '  +3.64%  '
Set-TextValue ($ws.Range('D39')) '18.20'
$ws.Range('E39').Value = '  +2.42%  '
Set-TextValue ($ws.Range('D40')) '0.305'
$ws.Range('E40').Value = '  +3.24%  '
Set-TextValue ($ws.Range('D41')) '1.55'
$ws.Range('E41').Value = '  +5.47%  '
Set-TextValue ($ws.Range('D42')) '4.47'
$ws.Range('E42').Value = '  +2.16%  '
Set-TextValue ($ws.Range('D43')) '1.10'
$ws.Range('E43').Value = '  +3.67%  '
Set-TextValue ($ws.Range('D44')) '2.15'
$ws.Range('E44').Value = '  +6.41%  '
Set-TextValue ($ws.Range('D45')) '3.40'
$ws.Range('E45').Value = '  +1.04%  '
Set-TextValue ($ws.Range('D46')) '132.56'
$ws.Range('E46').Value = '  +2.17%  '
Set-TextValue ($ws.Range('D47')) '0.0725'
$ws.Range('E47').Value = '  +0.49%  '
Set-TextValue ($ws.Range('D48')) '0.491'
$ws.Range('E48').Value = '  +2.20%  '
Set-TextValue ($ws.Range('D49')) '0.565'
$ws.Range('E49').Value = '  +0.54%  '
Set-TextValue ($ws.Range('D50')) '0.0931'
$ws.Range('E50').Value = '  +1.23%  '
$ws.Range('E51').Value = '  +2.80%  '
